$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") on rows 2-116 is a date column (serial date value).
# Update every existing value of 45189 (2023-09-20) to 45190 (2023-09-21).
$ws.Range("C2:C116").Value = 45190
